$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the slightly-adjusted timestamp on row 15 (A15)
$ws.Range("A15").Value = 45865.66692335648

# Append new row 16 with the latest sensor reading
$ws.Range("A16").Value = 45865.708647753
$ws.Range("A16").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("B16").Value = 2025
$ws.Range("C16").Value = 30
$ws.Range("D16").Value = 19
$ws.Range("E16").Value = 75.37
$ws.Range("F16").Value = 151.17
$ws.Range("G16").Value = 8.25
$ws.Range("H16").Value = "ESE"
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = "17:00:27"
